# "Moved the delete phase to before the conversion phase to make sure
# they don't tread on each other's toes" — the Omeka "Image" type label
# used on every map record is being renamed to the correct Dublin Core
# controlled-vocabulary term "Still Image".
#
# Deleting the stale "Image" shared-string value before writing the new
# "Still Image" text is what keeps Excel's shared-strings table tidy
# (old entry removed, new one appended) instead of ending up with both
# strings lingering side by side.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First Fleet Maps")

# Column A holds the Omeka item type for every record (header in A1,
# data in A2:A24). Replace the exact ("whole cell") value "Image" with
# "Still Image" everywhere it appears.
$rng = $ws.Range("A1:A24")
$rng.Replace("Image", "Still Image", $null, 1) | Out-Null

# Leave the sheet's selection where the editor ended up after making
# the change.
$ws.Activate()
$ws.Range("A29").Select()
